# Fix code systems and value sets for multi-modal communication
#
# * Metadata sheet: bump the "Date" value and the "Count" value (1 -> 8,
#   since the Concepts sheet now enumerates 8 concepts instead of 1).
# * Concepts sheet: replace the single placeholder concept row with the
#   full list of 8 concepts (7 new "Functional Communication Measure"
#   observation codes, plus the original "communicate-without-assistance"
#   concept which now becomes the final row).

$wb = $excel.ActiveWorkbook

$metaWs = $wb.Worksheets.Item("Metadata")
$conceptsWs = $wb.Worksheets.Item("Concepts")

# ---- Metadata sheet ---------------------------------------------------

# B8 = Date value. Not numeric-looking, so a plain string assignment is fine.
$metaWs.Cells.Item(8, 2).Value = "2022-04-05T11:13:11-04:00"

# B23 = Count value. A bare "8" would be auto-coerced to a Number by Excel's
# normal type inference, so a leading apostrophe forces text, matching the
# workbook's existing text-typed "Count" cell.
$metaWs.Cells.Item(23, 2).Value = "'8"

# ---- Concepts sheet -----------------------------------------------------
# Column layout: A=Level, B=Code, C=Display, D=Definition

# Give the new rows (3..9) the same formatting as the existing data row
# (row 2) before filling in any values, so borders/alignment match.
$conceptsWs.Range("A2:D2").Copy()
$conceptsWs.Range("A3:D9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @("1", "99829-4", "How often does the individual convey simple messages that are meaningful related to routine daily activities in LOW demand situations?", ""),
    @("1", "99830-2", "How often does the individual participate in short structured conversations that are meaningful in LOW demand situations?", ""),
    @("1", "99831-0", "How often does the individual convey complex messages that are meaningful in LOW demand situations?", ""),
    @("1", "99832-8", "How often does the individual convey simple messages that are meaningful related to routine daily activities in HIGH demand situations?", ""),
    @("1", "99833-6", "How often does the individual participate in short structured conversations that are meaningful in HIGH demand situations?", ""),
    @("1", "99834-4", "How often does the individual convey complex messages that are meaningful in HIGH demand situations?", ""),
    @("1", "99835-1", "Functional Communication Measure - Multi-Modal Functional Communication score [ASHA NOMS]", ""),
    @("1", "communicate-without-assistance", "Participate in communication exchanges without assistance", "How often does the individual participate in communication exchanges WITHOUT additional assistance from communication partner (no more than would be expected for chronological age)?")
)

$row = 2
foreach ($item in $data) {
    # Column A ("Level") is always the text "1". A bare "1" would be
    # auto-coerced to a Number, so force text with a leading apostrophe.
    $conceptsWs.Cells.Item($row, 1).Value = "'" + $item[0]
    $conceptsWs.Cells.Item($row, 2).Value = $item[1]
    $conceptsWs.Cells.Item($row, 3).Value = $item[2]

    if ($item[3] -ne "") {
        $conceptsWs.Cells.Item($row, 4).Value = $item[3]
    } else {
        $conceptsWs.Cells.Item($row, 4).ClearContents()
    }

    $row = $row + 1
}
